$d = $word.ActiveDocument

# Locate the paragraph that contains the misspelled "Javascript" mention
# that is flagged by the proofing-error markup (there is another,
# unrelated, "Javascript" substring elsewhere in the document, split
# across several runs, which must stay untouched).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Finalize the HTML search feature using Javascript*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$range = $target.Range
$fullText = $range.Text

$idx = $fullText.IndexOf("Javascript")
if ($idx -lt 0) {
    throw "Could not locate 'Javascript' inside target paragraph"
}
$beforeText = $fullText.Substring(0, $idx)
$afterWord = "JavaScript"

# Grab the paragraph's own opening tag (with all of its rsid/paraId
# attributes) plus its w:pPr so the rebuilt paragraph keeps its original
# formatting/identity - only the runs/text inside are changing.
$openXml = $range.WordOpenXML
if ($openXml -notmatch '(?s)(<w:p [^>]*>)') {
    throw "Could not extract paragraph opening tag"
}
$pOpenTag = $matches[1]

$pPr = ""
if ($openXml -match '(?s)(<w:pPr>.*?</w:pPr>)') {
    $pPr = $matches[1]
}

function Esc([string]$s) {
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# Rebuild the paragraph with exactly two runs: the untouched lead-in text
# (still carrying xml:space="preserve" since it has a trailing space) and
# a fresh run holding the corrected "JavaScript" - with no w:proofErr
# bookmarks around it, since the word is spelled correctly now.
$paraXml = $pOpenTag + $pPr +
    '<w:r><w:t xml:space="preserve">' + (Esc $beforeText) + '</w:t></w:r>' +
    '<w:r><w:t>' + (Esc $afterWord) + '</w:t></w:r>' +
    '</w:p>'

$package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$range.InsertXML($package)
